$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was added for "Zapallo italiano" (Femacal de La
# Calera). It belongs right after the current row 190, so insert a new row
# there; this pushes the old rows 191-219 down to 192-220, which is exactly
# what the diff shows (each old row's data now lives one row lower).
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new weekly record. The
# surrounding descriptive columns (A, B, C, E, F, G, H, I, N, O, Q, R) are
# identical to the record that is now at row 192, so only copy/set the
# columns that actually carry new data for this row.
$ws.Range("A191").Value = 3
$ws.Range("B191").Value = "Femacal de La Calera"
$ws.Range("C191").Value = "Coquimbo"
$ws.Range("D191").Value = "2021-10-05"
$ws.Range("E191").Value = 5
$ws.Range("F191").Value = 100112032
$ws.Range("G191").Value = "Zapallo italiano"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 85
$ws.Range("K191").Value = 14000
$ws.Range("L191").Value = 15000
$ws.Range("M191").Value = 14529
$ws.Range("N191").Value = "$/caja 70 unidades"
$ws.Range("O191").Value = "Región de Arica y Parinacota"
$ws.Range("P191").Value = 208
$ws.Range("Q191").Value = 70
$ws.Range("R191").Value = "Hortaliza"
